$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F2, F5, F9, F10 increment
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 337
$ws1.Range("F5").Value = 4968
$ws1.Range("F9").Value = 754
$ws1.Range("F10").Value = 237

# Sheet "全部类型" (sheet4): F2, F5, F9, F11 increment
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 337
$ws4.Range("F5").Value = 4968
$ws4.Range("F9").Value = 754
$ws4.Range("F11").Value = 237
